$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44230
$ws.Range("J2").Value = 250

# Row 3
$ws.Range("D3").Value = 44875
$ws.Range("J3").Value = 90
$ws.Range("K3").Value = 7000
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = 7000
$ws.Range("P3").Value = 438

# Row 4
$ws.Range("D4").Value = 44189
$ws.Range("J4").Value = 250

# Row 6
$ws.Range("D6").Value = 44187
$ws.Range("J6").Value = 160

# Row 7
$ws.Range("D7").Value = 44210
$ws.Range("J7").Value = 340
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = 5500
$ws.Range("P7").Value = 344

# Row 8
$ws.Range("D8").Value = 44236
$ws.Range("J8").Value = 180
$ws.Range("K8").Value = 4000
$ws.Range("L8").Value = 4500
$ws.Range("M8").Value = 4167
$ws.Range("O8").Value = "Región Metropolitana"
$ws.Range("P8").Value = 260

# Row 9
$ws.Range("D9").Value = 44873
$ws.Range("K9").Value = 8000
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = 8000
$ws.Range("P9").Value = 500

# Row 10
$ws.Range("D10").Value = 44188
$ws.Range("J10").Value = 210
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 6000
$ws.Range("M10").Value = 5500
$ws.Range("O10").Value = "Provincia de Quillota"
$ws.Range("P10").Value = 344

# Row 11
$ws.Range("D11").Value = 44208
$ws.Range("J11").Value = 160
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = 5500
$ws.Range("O11").Value = "Provincia de Quillota"
$ws.Range("P11").Value = 344

# Row 12
$ws.Range("D12").Value = 44231
$ws.Range("L12").Value = 6000
$ws.Range("M12").Value = 5500
$ws.Range("P12").Value = 344

# Row 13
$ws.Range("D13").Value = 44232
$ws.Range("J13").Value = 250
$ws.Range("K13").Value = 5000
$ws.Range("M13").Value = 5500
$ws.Range("O13").Value = "Provincia de Quillota"
$ws.Range("P13").Value = 344

# Row 14
$ws.Range("D14").Value = 44883
$ws.Range("J14").Value = 180
$ws.Range("K14").Value = 7000
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = 7500
$ws.Range("P14").Value = 469

# Row 15
$ws.Range("D15").Value = 44292
$ws.Range("J15").Value = 90
$ws.Range("K15").Value = 6000
$ws.Range("M15").Value = 6000
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 375

# Row 16
$ws.Range("D16").Value = 44204
$ws.Range("J16").Value = 430
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 6000
$ws.Range("M16").Value = 5500
$ws.Range("P16").Value = 344

# Row 17
$ws.Range("D17").Value = 44846
$ws.Range("L17").Value = 5000
$ws.Range("M17").Value = 5000
$ws.Range("P17").Value = 312

# Row 18
$ws.Range("D18").Value = 44882
$ws.Range("J18").Value = 70
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 7000
$ws.Range("P18").Value = 438

# Row 19
$ws.Range("D19").Value = 44215
$ws.Range("K19").Value = 5000
$ws.Range("L19").Value = 6000
$ws.Range("M19").Value = 5500
$ws.Range("P19").Value = 344

# Row 20
$ws.Range("D20").Value = 44855
$ws.Range("J20").Value = 70
$ws.Range("K20").Value = 6000
$ws.Range("L20").Value = 7000
$ws.Range("M20").Value = 6500
$ws.Range("P20").Value = 406

# Row 21
$ws.Range("D21").Value = 44251
$ws.Range("J21").Value = 120
$ws.Range("L21").Value = 5000
$ws.Range("M21").Value = 5000
$ws.Range("O21").Value = "Región Metropolitana"
$ws.Range("P21").Value = 312

